$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Liz's row (row 2; row 1 is the header row) ---
$t.Cell(2, 2).Range.Text = "Finished User and their profiles. Created a profile page"
$t.Cell(2, 3).Range.Text = "Work on the inventory page"

$obstacleCell = $t.Cell(2, 4)
$obstacleCell.Range.Text = "none"
# Word drops a "_GoBack" bookmark at the site of the most recent edit.
$d.Bookmarks.Add("_GoBack", $obstacleCell.Range)

# --- Ian's row (row 3): the "what did you do" text got retyped/edited
# mid-sentence, leaving the sentence split across two runs with identical
# formatting but the same overall text. ---
$ianCell = $t.Cell(3, 2)
$cellStart = $ianCell.Range.Start
$part1 = "Compiled a short list of tools that I "
$part2 = "found images and descriptions for."
$ianCell.Range.Text = $part1 + $part2

$splitAt = $cellStart + $part1.Length
$secondRun = $d.Range($splitAt, $splitAt + $part2.Length)
# Toggling a character property and back forces the run to split at this
# boundary without altering the visible formatting.
$secondRun.Bold = 1
$secondRun.Bold = 0
